# Update the author list on the cover page from
# "Daniel Ricard and Peter Comeau" to the full author list, with each
# author/separator rendered as its own run (mirroring how Word leaves
# behind a trail of individually-edited runs), while preserving the
# existing run formatting (Microsoft Sans Serif, sz 28, noProof).

$d = $word.ActiveDocument

# Locate the existing author run and replace its text with the first
# author's chunk; this also collapses/repositions $rng to track it.
$rng = $d.Content
$rng.Find.Execute("Daniel Ricard and Peter Comeau")
$rng.Text = "Daniel Ricard,"
$rng.Collapse(0)

# Remaining chunks, inserted in order right after the first one. Each
# chunk ends up adjacent to its neighbour, so Word's live run-coalescing
# will merge them all into a single run for now -- that's fixed up in
# the second pass below.
$chunks = @(
    " Aaron, Adamack, ",
    "Peter Comeau,",
    " ",
    "Allan Debertin,",
    " ",
    "Kim Emond,",
    " ",
    "Tracey Loewen,",
    " ",
    "Gregory Puncher,",
    " ",
    "Meredith Schofield,",
    " ",
    "Andrew Smith and",
    " ",
    "Stephen Wischniowski"
)

# Track the document-character offset of the boundary after each chunk
# so we can revisit each span once all the text is in place.
$offsets = @($rng.Start)

foreach ($chunk in $chunks) {
    $rng.InsertAfter($chunk)
    $offsets += $rng.End
    $rng.Collapse(0)
}

# Second pass: re-select each inserted chunk (now contiguous, identically
# formatted text) and toggle Bold on/off. The toggle is a no-op on the
# rendered formatting, but it forces the run to stay distinct from its
# neighbours instead of being silently coalesced back into one run --
# reproducing the many separate <w:r> elements from the target revision.
for ($i = 0; $i -lt $offsets.Length - 1; $i++) {
    $seg = $d.Range($offsets[$i], $offsets[$i + 1])
    $seg.Font.Bold = 1
    $seg.Font.Bold = 0
}
